# Update countries & provincias Spain
# Applies the 17-Sep-2020 data refresh:
#  - three country-name rank swaps (Malasia/Namibia, Surinam/Ruanda,
#    Islas Malvinas/Montserrat) caused by the underlying data being
#    re-sorted by "Casos totales" after the refresh
#  - the "Datos actualizados" timestamp text
#  - refreshed numeric figures for the affected country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 17 de Septiembre de 2020 a las 22:40"

# --- Country name swaps (rank changed after the refresh) ---
# Malasia <-> Namibia
$ws.Cells.Item(97,1).Value  = "Namibia"
$ws.Cells.Item(98,1).Value  = "Malasia"

# Surinam <-> Ruanda
$ws.Cells.Item(123,1).Value = "Ruanda"
$ws.Cells.Item(124,1).Value = "Surinam"

# Islas Malvinas <-> Montserrat
$ws.Cells.Item(214,1).Value = "Montserrat"
$ws.Cells.Item(215,1).Value = "Islas Malvinas"

# --- Refreshed numeric values ---
# Estados Unidos (row 4)
$ws.Cells.Item(4,2).Value = 6862145
$ws.Cells.Item(4,3).Value = 33844
$ws.Cells.Item(4,4).Value = 4141167
$ws.Cells.Item(4,5).Value = 2518998
$ws.Cells.Item(4,7).Value = 632
$ws.Cells.Item(4,8).Value = 201980

# India (row 5)
$ws.Cells.Item(5,2).Value = 5212686
$ws.Cells.Item(5,3).Value = 96793
$ws.Cells.Item(5,4).Value = 4109828
$ws.Cells.Item(5,5).Value = 1018454
$ws.Cells.Item(5,7).Value = 1174
$ws.Cells.Item(5,8).Value = 84404

# Francia (row 15)
$ws.Cells.Item(15,4).Value = 90840
$ws.Cells.Item(15,5).Value = 293546

# Alemania (row 25)
$ws.Cells.Item(25,2).Value = 269035
$ws.Cells.Item(25,3).Value = 2170
$ws.Cells.Item(25,5).Value = 20480
$ws.Cells.Item(25,7).Value = 6
$ws.Cells.Item(25,8).Value = 9455

# Israel (row 27)
$ws.Cells.Item(27,2).Value = 175256
$ws.Cells.Item(27,3).Value = 4791
$ws.Cells.Item(27,4).Value = 126329
$ws.Cells.Item(27,5).Value = 47758
$ws.Cells.Item(27,7).Value = 8
$ws.Cells.Item(27,8).Value = 1169

# Costa de Marfil (row 83)
$ws.Cells.Item(83,2).Value = 19158
$ws.Cells.Item(83,3).Value = 26
$ws.Cells.Item(83,4).Value = 18330
$ws.Cells.Item(83,5).Value = 708

# Namibia (row 97, after swap)
$ws.Cells.Item(97,2).Value = 10078
$ws.Cells.Item(97,3).Value = 114
$ws.Cells.Item(97,4).Value = 7685
$ws.Cells.Item(97,5).Value = 2285
$ws.Cells.Item(97,8).Value = 108

# Malasia (row 98, after swap)
$ws.Cells.Item(98,2).Value = 10052
$ws.Cells.Item(98,3).Value = 21
$ws.Cells.Item(98,4).Value = 9250
$ws.Cells.Item(98,5).Value = 674
$ws.Cells.Item(98,8).Value = 128

# Ruanda (row 123, after swap)
$ws.Cells.Item(123,2).Value = 4653
$ws.Cells.Item(123,3).Value = 19
$ws.Cells.Item(123,4).Value = 2817
$ws.Cells.Item(123,5).Value = 1813
$ws.Cells.Item(123,7).Value = 1
$ws.Cells.Item(123,8).Value = 23

# Surinam (row 124, after swap)
$ws.Cells.Item(124,2).Value = 4645
$ws.Cells.Item(124,4).Value = 4089
$ws.Cells.Item(124,5).Value = 461
$ws.Cells.Item(124,8).Value = 95

# Angola (row 129)
$ws.Cells.Item(129,2).Value = 3789
$ws.Cells.Item(129,3).Value = 114
$ws.Cells.Item(129,4).Value = 1405
$ws.Cells.Item(129,5).Value = 2240
$ws.Cells.Item(129,7).Value = 1
$ws.Cells.Item(129,8).Value = 144

# Siria (row 130)
$ws.Cells.Item(130,2).Value = 3691
$ws.Cells.Item(130,3).Value = 37
$ws.Cells.Item(130,4).Value = 903
$ws.Cells.Item(130,5).Value = 2623
$ws.Cells.Item(130,7).Value = 2
$ws.Cells.Item(130,8).Value = 165

# Montserrat (row 214, after swap)
$ws.Cells.Item(214,4).Value = 12
$ws.Cells.Item(214,8).Value = 1

# Islas Malvinas (row 215, after swap)
$ws.Cells.Item(215,4).Value = 13
$ws.Cells.Item(215,8).Value = 0
